# Generate Report for Handback
# Update timestamp cells that track handoff/handback generation times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for first row
$wsOverview.Range("G2").Value = "2016-08-17 04:58:55"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first row
$wsZhCn.Range("H2").Value = "2016-08-17 04:58:50"
$wsZhCn.Range("K2").Value = "2016-08-17 04:59:12"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first row
$wsDeDe.Range("H2").Value = "2016-08-17 04:58:55"
$wsDeDe.Range("K2").Value = "2016-08-17 04:59:20"
